$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) for the HP 国风动漫游戏嘉年华 row (F4)
# and the 动漫游戏展 row (F5) on both the "展览" and "全部类型" sheets.
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 164
    $ws.Range("F5").Value = 11
}
